$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vocab")

# --- Column A (German): entered first, in row order 283-290 ---
$ws.Cells.Item(283, 1).Value = '…spricht uns besonders an'
$ws.Cells.Item(284, 1).Value = '…haben wir / habe ich die Möglichkeit…'
$rt284 = $ws.Cells.Item(284, 1).Characters(28, 11)
$rt284.Font.Name = "Calibri"
$rt284.Font.Size = 11
$ws.Cells.Item(285, 1).Value = 'Für…wünschen wir uns…'
$rt285 = $ws.Cells.Item(285, 1).Characters(2, 20)
$rt285.Font.Name = "Calibri"
$rt285.Font.Size = 11
$ws.Cells.Item(286, 1).Value = 'Ich / Wir können uns vorstellen, …'
$rt286 = $ws.Cells.Item(286, 1).Characters(12, 23)
$rt286.Font.Name = "Calibri"
$rt286.Font.Size = 11
$ws.Cells.Item(287, 1).Value = 'Wir hätten gerne…'
$rt287 = $ws.Cells.Item(287, 1).Characters(6, 12)
$rt287.Font.Name = "Calibri"
$rt287.Font.Size = 11
$ws.Cells.Item(288, 1).Value = 'Auβerdem…'
$rt288 = $ws.Cells.Item(288, 1).Characters(3, 7)
$rt288.Font.Name = "Calibri"
$rt288.Font.Size = 11
$ws.Cells.Item(289, 1).Value = 'Ich/Wir finden… gut, weil/da…'
$ws.Cells.Item(290, 1).Value = 'Besonders wichtig ist mir/uns…'

# --- Column B (English): entered after all German phrases ---
$ws.Cells.Item(283, 2).Value = 'appeals to me / us'
$ws.Cells.Item(284, 2).Value = 'we have the possibility'
$ws.Cells.Item(285, 2).Value = 'for… we wish'
$ws.Cells.Item(286, 2).Value = 'we can imagine / picture'
$ws.Cells.Item(287, 2).Value = 'we would like…'
$ws.Cells.Item(288, 2).Value = 'besides that…'
$ws.Cells.Item(289, 2).Value = 'I find…good, because'
$ws.Cells.Item(290, 2).Value = 'It is particularly important to us…'

# --- Columns C/D/E (lesson metadata) ---
$ws.Cells.Item(283, 3).Value = 44657
$ws.Cells.Item(283, 4).Value = 10
$ws.Cells.Item(283, 5).Value = "phrase"
$ws.Cells.Item(284, 3).Value = 44657
$ws.Cells.Item(284, 4).Value = 10
$ws.Cells.Item(284, 5).Value = "phrase"
$ws.Cells.Item(285, 3).Value = 44657
$ws.Cells.Item(285, 4).Value = 10
$ws.Cells.Item(285, 5).Value = "phrase"
$ws.Cells.Item(286, 3).Value = 44657
$ws.Cells.Item(286, 4).Value = 10
$ws.Cells.Item(286, 5).Value = "phrase"
$ws.Cells.Item(287, 3).Value = 44657
$ws.Cells.Item(287, 4).Value = 10
$ws.Cells.Item(287, 5).Value = "phrase"
$ws.Cells.Item(288, 3).Value = 44657
$ws.Cells.Item(288, 4).Value = 10
$ws.Cells.Item(288, 5).Value = "phrase"
$ws.Cells.Item(289, 3).Value = 44657
$ws.Cells.Item(289, 4).Value = 10
$ws.Cells.Item(289, 5).Value = "phrase"
$ws.Cells.Item(290, 3).Value = 44657
$ws.Cells.Item(290, 4).Value = 10
$ws.Cells.Item(290, 5).Value = "phrase"

# Copy the date style (style index reused, matching numFmtId 14 m/d/yyyy) from C282 onto the new date cells
$ws.Range("C282").Copy()
$ws.Range("C283:C290").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update view/selection to match the saved state
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 241
$win.ScrollColumn = 1
$ws.Range("C282:E290").Select() | Out-Null
